$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the target cells are treated as text so Excel does not
# auto-convert number-like strings (e.g. "208.48") into numeric values.
$updateRange = $ws.Range("D2:E51")
$updateRange.NumberFormat = "@"

$ws.Range("D2").Value = '26.155.72'
$ws.Range("E2").Value = '  -2.21%  '
$ws.Range("D3").Value = '1.574.53'
$ws.Range("E3").Value = '  -1.80%  '
$ws.Range("E4").Value = '  -0.52%  '
$ws.Range("D5").Value = '208.48'
$ws.Range("E5").Value = '  -1.67%  '
$ws.Range("E6").Value = '  -2.89%  '
$ws.Range("E7").Value = '  -0.49%  '
$ws.Range("E8").Value = '  -1.74%  '
$ws.Range("E9").Value = '  -1.37%  '
$ws.Range("D10").Value = '19.58'
$ws.Range("E10").Value = '  -0.69%  '
$ws.Range("E11").Value = '  -0.33%  '
$ws.Range("D12").Value = '1.794.61'
$ws.Range("E12").Value = '  -1.88%  '
$ws.Range("D13").Value = '1.591.43'
$ws.Range("E13").Value = '  -0.84%  '
$ws.Range("D14").Value = '4.06'
$ws.Range("E14").Value = '  -0.48%  '
$ws.Range("D15").Value = '0.514'
$ws.Range("E15").Value = '  -2.23%  '
$ws.Range("D16").Value = '64.33'
$ws.Range("E16").Value = '  -1.13%  '
$ws.Range("D17").Value = '26.140.40'
$ws.Range("E17").Value = '  -2.20%  '
$ws.Range("E18").Value = '  -2.49%  '
$ws.Range("D19").Value = '7.27'
$ws.Range("E19").Value = '  +1.84%  '
$ws.Range("D20").Value = '207.70'
$ws.Range("E20").Value = '  -1.07%  '
$ws.Range("E21").Value = '  -0.44%  '
$ws.Range("D22").Value = '4.25'
$ws.Range("E22").Value = '  -1.36%  '
$ws.Range("E23").Value = '  -2.83%  '
$ws.Range("D24").Value = '8.83'
$ws.Range("E24").Value = '  -2.77%  '
$ws.Range("D25").Value = '143.46'
$ws.Range("E25").Value = '  -0.28%  '
$ws.Range("E26").Value = '  -0.43%  '
$ws.Range("D27").Value = '6.98'
$ws.Range("E27").Value = '  -1.74%  '
$ws.Range("E28").Value = '  -1.88%  '
$ws.Range("D29").Value = '15.21'
$ws.Range("E29").Value = '  -1.03%  '
$ws.Range("E30").Value = '  -0.69%  '
$ws.Range("E31").Value = '  -1.46%  '
$ws.Range("E32").Value = '  -2.03%  '
$ws.Range("D33").Value = '2.97'
$ws.Range("E33").Value = '  +0.18%  '
$ws.Range("D34").Value = '1.275.64'
$ws.Range("E34").Value = '  -1.22%  '
$ws.Range("D35").Value = '0.614'
$ws.Range("E35").Value = '  +3.17%  '
$ws.Range("E36").Value = '  -1.68%  '
$ws.Range("E37").Value = '  -0.98%  '
$ws.Range("D38").Value = '0.0166'
$ws.Range("E38").Value = '  -2.58%  '
$ws.Range("E39").Value = '  -11.17%  '
$ws.Range("E40").Value = '  -2.61%  '
$ws.Range("D41").Value = '5.55'
$ws.Range("E41").Value = '  +1.95%  '
$ws.Range("E42").Value = '  -2.78%  '
$ws.Range("D43").Value = '0.763'
$ws.Range("E43").Value = '  -2.29%  '
$ws.Range("D44").Value = '62.06'
$ws.Range("E44").Value = '  -1.72%  '
$ws.Range("D45").Value = '1.708.25'
$ws.Range("E45").Value = '  -1.91%  '
$ws.Range("D46").Value = '88.90'
$ws.Range("E46").Value = '  -1.75%  '
$ws.Range("E47").Value = '  +0.61%  '
$ws.Range("E48").Value = '  -3.12%  '
$ws.Range("E49").Value = '  -2.21%  '
$ws.Range("E50").Value = '  -1.59%  '
$ws.Range("E51").Value = '  -0.45%  '

# Restore the default (Normal) style so no stray number-format override
# is left behind on these cells.
$updateRange.Style = "Normal"
